$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataValueSet")

# Add the new data row (row 2) with payload values for dataElement/optionSet/options post method
$ws.Range("A2").Value = "pGeBz8X2jRq"
$ws.Range("B2").Value = "yHSAPCLxecr"
$ws.Range("C2").Value = "WTSe3FmRFmD"
$ws.Range("D2").Value = 202007
$ws.Range("E2").Value = 8

# Update selection to match the post-edit active cell state
$ws.Range("A2").Select()

$wb.Save()
